$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated simulation results for the case with 380 kV (row 2-25, columns C-K, O)
$data = @{
    2 = @{ "C"=0.02576265356624674; "D"=0.2187163449224272; "E"=0.1739981335038863; "F"=1.244428632407768; "G"=0.002447845095296664; "I"=0.6196914913318281; "J"=0.1829694901036945; "K"=1.926523897119182; "O"=2.881358039394911 }
    3 = @{ "C"=0.02308503200887202; "D"=0.2111408010845963; "E"=0.1693529599698778; "F"=1.254511569362656; "G"=0.002450909747672517; "I"=0.6263383624242458; "J"=0.1791880082050881; "K"=1.709701206097748; "O"=2.922835161227283 }
    4 = @{ "C"=0.02143188334719781; "D"=0.2065450991217261; "E"=0.1665791796037865; "F"=1.261676498512166; "G"=0.002452889805782892; "I"=0.6309749586581397; "J"=0.1769782326122424; "K"=1.576197383112344; "O"=2.950833808790989 }
    5 = @{ "C"=0.02075596059319196; "D"=0.2046864526286925; "E"=0.1654685898885475; "F"=1.264840722367524; "G"=0.002453721500704383; "I"=0.6330037209054247; "J"=0.1761058865288092; "K"=1.521702919598681; "O"=2.962878484742234 }
    6 = @{ "C"=0.0206435891021215; "D"=0.2043786836306936; "E"=0.1652853712642717; "F"=1.265380890611013; "G"=0.002453861103332126; "I"=0.6333489992642072; "J"=0.1759627339152132; "K"=1.512648786959346; "O"=2.964916800276711 }
    7 = @{ "C"=0.02142277669020842; "D"=0.206519975341493; "E"=0.1665641217878253; "F"=1.261718183151203; "G"=0.002452900921651429; "I"=0.6310017556489491; "J"=0.1769663538807293; "K"=1.575462813862373; "O"=2.950993678512702 }
    8 = @{ "C"=0.024841311841449; "D"=0.2160928092285843; "E"=0.1723802331284645; "F"=1.247702808308652; "G"=0.002448881420026642; "I"=0.6218678531275792; "J"=0.181642358062625; "K"=1.851842994213996; "O"=2.895132988876298 }
    9 = @{ "C"=0.03147195897220456; "D"=0.2353025653250995; "E"=0.1844065009046929; "F"=1.227967899634884; "G"=0.002441776126493937; "I"=0.608378883696183; "J"=0.1917030994898568; "K"=2.390734903853456; "O"=2.805747503124763 }
    10 = @{ "C"=0.03629799706533277; "D"=0.2496781250933822; "E"=0.1936205659822647; "F"=1.218223213511692; "G"=0.002437024697682038; "I"=0.6011868586913991; "J"=0.1996416610974592; "K"=2.784651785247092; "O"=2.752461608117869 }
    11 = @{ "C"=0.0384834437947319; "D"=0.2562740234660907; "E"=0.1978944937507165; "F"=1.214829124931896; "G"=0.002434963925697048; "I"=0.598510065754482; "J"=0.2033727440979476; "K"=2.963395938758708; "O"=2.730930803019049 }
    12 = @{ "C"=0.03930956221374515; "D"=0.2587797229401474; "E"=0.199524748364432; "F"=1.21369377470451; "G"=0.002434197964931976; "I"=0.5975823638502291; "J"=0.2048028818675647; "K"=3.031014169234368; "O"=2.723168957821542 }
    13 = @{ "C"=0.03913170841413205; "D"=0.2582397230833919; "E"=0.1991731191791573; "F"=1.213931616534182; "G"=0.002434362288464831; "I"=0.5977783324827683; "J"=0.2044941080658162; "K"=3.016454466415439; "O"=2.724823173308152 }
    14 = @{ "C"=0.0385514385592387; "D"=0.2564800097725879; "E"=0.1980283792894468; "F"=1.214732711345732; "G"=0.00243490062130176; "I"=0.5984320184987482; "J"=0.2034900563419768; "K"=2.968960320021665; "O"=2.730284376611678 }
    15 = @{ "C"=0.03819581503911706; "D"=0.2554031693610455; "E"=0.197328730048298; "F"=1.215242944638661; "G"=0.002435232240089734; "I"=0.5988436239213399; "J"=0.2028772940666528; "K"=2.93985978604087; "O"=2.733680547619144 }
    16 = @{ "C"=0.03615497007419322; "D"=0.2492481908957558; "E"=0.1933429096078854; "F"=1.21846597475664; "G"=0.002437161393816075; "I"=0.6013737981148637; "J"=0.1994002392681438; "K"=2.772961066304617; "O"=2.753923368975052 }
    17 = @{ "C"=0.03490040762503099; "D"=0.2454866579076906; "E"=0.1909188157292476; "F"=1.220709647739255; "G"=0.002438370604133397; "I"=0.603078627245246; "J"=0.197297882755521; "K"=2.670456185859621; "O"=2.767036876189053 }
    18 = @{ "C"=0.03417788263195121; "D"=0.2433284385756593; "E"=0.1895322995144184; "F"=1.222097895532627; "G"=0.002439075590346782; "I"=0.604115170597602; "J"=0.196099936328892; "K"=2.611455882269979; "O"=2.77483431414592 }
    19 = @{ "C"=0.03393308871207523; "D"=0.2425986194868983; "E"=0.1890641826556134; "F"=1.222584704063628; "G"=0.002439315916801953; "I"=0.6044757273064647; "J"=0.1956962675117353; "K"=2.591472251658161; "O"=2.77751811567498 }
    20 = @{ "C"=0.03503405483219524; "D"=0.2458865303996163; "E"=0.191176061970296; "F"=1.220460684635981; "G"=0.002438240901141098; "I"=0.6028913496380639; "J"=0.1975205152087085; "K"=2.681372399238285; "O"=2.765614527023217 }
    21 = @{ "C"=0.03872191784513745; "D"=0.256996664810373; "E"=0.1983642969037547; "F"=1.21449333739136; "G"=0.002434742109474508; "I"=0.5982376794893582; "J"=0.2037845020511639; "K"=2.982912383437849; "O"=2.72866965173813 }
    22 = @{ "C"=0.04112360796951009; "D"=0.2643042092662995; "E"=0.2031310495434013; "F"=1.211467413092805; "G"=0.0024325394004392; "I"=0.5956973530405207; "J"=0.2079789850429989; "K"=3.17958649291279; "O"=2.70680628235948 }
    23 = @{ "C"=0.03984257313963724; "D"=0.2603998302590185; "E"=0.2005806590145909; "F"=1.213002250186079; "G"=0.002433707368562779; "I"=0.5970071902352529; "J"=0.2057310945679376; "K"=3.074655536487683; "O"=2.718265747516767 }
    24 = @{ "C"=0.03497363685374921; "D"=0.245705734532379; "E"=0.1910597387412807; "F"=1.220572934597683; "G"=0.002438299509317458; "I"=0.6029758421440476; "J"=0.1974198296688172; "K"=2.676437393474146; "O"=2.766256766587702 }
    25 = @{ "C"=0.02968611704844193; "D"=0.2300594583771982; "E"=0.1810866232795689; "F"=1.232473961184752; "G"=0.002443615619215539; "I"=0.6115522416078818; "J"=0.1888856223532258; "K"=2.245294057863532; "O"=2.827760619217344 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Range("$col$row").Value = $data[$row][$col]
    }
}
